# StagingTemplates/Staging.Organization.xlsx update:
# "shortname longname staging template update for Organization"
#
# The Organization staging template's header row is changed from:
#   Organization_ID | BusinessKey | OrganizationTypeBusinessKey | Code | Name | ParentOrganization_ID
# to:
#   Organization_ID | BusinessKey | OrganizationTypeBusinessKey | Code | ShortName | LongName | ParentOrganization_ID
#
# i.e. the single "Name" column is replaced by two columns ("ShortName" and
# "LongName"), and "ParentOrganization_ID" is pushed out one column to the
# right (from F to G) to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "ParentOrganization_ID" header from F2 to the new G2
# column (shifting it one column to the right).
$ws.Range("G2").Value = $ws.Range("F2").Text

# Replace the old "Name" column (E2) and repurpose the vacated F2 cell with
# the two new headers.
$ws.Range("E2").Value = "ShortName"
$ws.Range("F2").Value = "LongName"

# Match the saved selection/active cell shown in the workbook (E2).
$ws.Range("E2").Select() | Out-Null

# Resize the columns so the headers are fully visible (best-fit), matching
# the widths Excel computed for the new header text.
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 27.333333333333332
$ws.Columns.Item(4).ColumnWidth = 4.666666666666667
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666
$ws.Columns.Item(6).ColumnWidth = 9.666666666666666
$ws.Columns.Item(7).ColumnWidth = 20.666666666666668
